# Generate Report for handback
# The ce536fa2-b3da-4f14-a74d-08faf91cb8e1.md file has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# status + the "Latest Handback DateTime" columns on the per-locale sheets,
# and roll the status up onto the Overview sheet too.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-01-25 10:57:22"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-01-25 10:57:37"
